$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 17
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 31
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0
$ws.Range("B9").Value = 11
$ws.Range("E9").Value = 1
$ws.Range("B10").Value = 16
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("B11").Value = 8
$ws.Range("C11").Value = 19
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 4
$ws.Range("B12").Value = 9
$ws.Range("C12").Value = 39
$ws.Range("B13").Value = 8
$ws.Range("C13").Value = 24
$ws.Range("D13").Value = 2
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 35
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 20
$ws.Range("B17").Value = 8
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 3
$ws.Range("C18").Value = 22
$ws.Range("B19").Value = 8
$ws.Range("C19").Value = 8
$ws.Range("E19").Value = 0
$ws.Range("B20").Value = 3
$ws.Range("C20").Value = 5
$ws.Range("B21").Value = 7
$ws.Range("F23").Value = 0
$ws.Range("B24").Value = 9
$ws.Range("F24").Value = 3
$ws.Range("B25").Value = 9
$ws.Range("C25").Value = 37
$ws.Range("B26").Value = 6
$ws.Range("C26").Value = 36
$ws.Range("B28").Value = 6
$ws.Range("E28").Value = 1
$ws.Range("F31").Value = 3
